# TODO maj 3eme miaritory
# Adds a new worksheet "3 eme fieretatory - contrat" after "2eme fieretatory",
# fills it with the contrat-essai TODO rows, and makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the last existing sheet ("2eme fieretatory")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "3 eme fieretatory - contrat"

# Fill column A rows 3-11 first, then column B rows 3-11 (matches authoring order
# so new shared-string entries land in the same order as the source edit)
$ws.Cells.Item(3, 1).Value = "liste des entretiens"
$ws.Cells.Item(4, 1).Value = "detail entretien"
$ws.Cells.Item(5, 1).Value = "choix d'embauche"
$ws.Cells.Item(6, 1).Value = "liste futur emp"
$ws.Cells.Item(7, 1).Value = "creation contrat essai"
$ws.Cells.Item(8, 1).Value = "resume contrat essai"
$ws.Cells.Item(9, 1).Value = "liste contrat essai"
$ws.Cells.Item(10, 1).Value = "detail contrat essai"
$ws.Cells.Item(11, 1).Value = "chgt contrat - choix contrat"

$ws.Cells.Item(3, 2).Value = "ny avo"
$ws.Cells.Item(4, 2).Value = "ny avo"
$ws.Cells.Item(5, 2).Value = "harena"
$ws.Cells.Item(6, 2).Value = "harena"
$ws.Cells.Item(7, 2).Value = "harena"
$ws.Cells.Item(8, 2).Value = "harena"
$ws.Cells.Item(10, 2).Value = "ny avo"
$ws.Cells.Item(11, 2).Value = "ny avo"

# Header row
$ws.Cells.Item(1, 1).Value = "Taches"
$ws.Cells.Item(1, 2).Value = "Qui"

# Row 2
$ws.Cells.Item(2, 1).Value = "maj insertion besoin recrutement"
$ws.Cells.Item(2, 2).Value = "ny avo"

# B9 carries a trailing space, distinct shared-string entry from "ny avo"
$ws.Cells.Item(9, 2).Value = "ny avo "

# Column A width
$ws.Columns.Item(1).ColumnWidth = 33.66

# Selection left on the new sheet
$ws.Range("C10").Select() | Out-Null
